# Generate Report for Archive
#
# The localization status for the zh-cn and de-de targets moved from
# "Ready for handoff" to "In Translation". This text lives in a shared
# string that is referenced from:
#   - Overview!E2 (zh-cn status) and Overview!F2 (de-de status)
#   - zh-cn!C2 (Status column)
#   - de-de!C2 (Status column)
# Updating the text causes the "Status" columns to re-autofit narrower
# (the new text is shorter than the old text).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# Update the status values everywhere they appear.
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsZhCn.Range("C2").Value = $newStatus
$wsDeDe.Range("C2").Value = $newStatus

# Re-fit the Status columns to the new (shorter) content.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
